# "Adicionados balanços concatenados em uma única planilha."
# Extend the NINJ3 balance-sheet table with three more reporting periods:
# column N = 31/12/2023, column O = 31/03/2024, column P = 30/06/2024.
# The sheet's used range grows from A1:M80 to A1:P80.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new period labels, styled like the existing
#     header cells (bold / bordered / centered, same format as M1). ---
$ws.Range("M1").Copy()
$ws.Range("N1:P1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("N1").Value = "31/12/2023"
$ws.Range("O1").Value = "31/03/2024"
$ws.Range("P1").Value = "30/06/2024"

# --- Data rows: N/O/P values for each line item (row number -> N,O,P). ---
$data = @{
    2  = @(297688.992, 301051.008, 307979.008)
    3  = @(287040.992, 291636, 299833.984)
    4  = @(24285, 20818, 5876)
    5  = @(251712, 256988, 278312.992)
    6  = @(4104, 5065, 4464)
    7  = @(0, 0, 0)
    8  = @(0, 0, 0)
    9  = @(0, 0, 0)
    10 = @(0, 0, 0)
    11 = @(6940, 8765, 11181)
    12 = @(0, 0, 0)
    13 = @(0, 0, 0)
    14 = @(0, 0, 0)
    15 = @(0, 0, 0)
    16 = @(0, 0, 0)
    17 = @(0, 0, 0)
    18 = @(0, 0, 0)
    19 = @(0, 0, 0)
    20 = @(0, 0, 0)
    21 = @(0, 0, 0)
    22 = @(0, 0, 0)
    23 = @(2825, 2604, 2347)
    24 = @(7823, 6811, 5798)
    25 = @(0, 0, 0)
    26 = @(297688.992, 301051.008, 307979.008)
    27 = @(16213, 16680, 17272)
    28 = @(0, 0, 0)
    29 = @(5700, 6710, 6264)
    30 = @(0, 0, 0)
    31 = @(0, 0, 0)
    32 = @(0, 0, 0)
    33 = @(0, 0, 0)
    34 = @(10513, 9970, 11008)
    35 = @(0, 0, 0)
    36 = @(0, 0, 0)
    37 = @(1141, 1048, 878)
    38 = @(1042, 949, 832)
    39 = @(0, 0, 0)
    40 = @(0, 0, 0)
    41 = @(0, 0, 0)
    42 = @(0, 0, 0)
    43 = @(99, 99, 46)
    44 = @(0, 0, 0)
    45 = @(0, 0, 0)
    46 = @(0, 0, 0)
    47 = @(280335.008, 283323.008, 289828.992)
    48 = @(267386, 267386, 267387.008)
    49 = @(8341, 8423, 8377)
    50 = @(0, 0, 0)
    51 = @(0, 0, 0)
    52 = @(4608, 7514, 14065)
    53 = @(0, 0, 0)
    54 = @(0, 0, 0)
    55 = @(0, 0, 0)
    56 = @(0, 0, 0)
    59 = @(14424, 16684, 15742)
    60 = @(-2893, -2499, -2463)
    61 = @(11531, 14185, 13279)
    62 = @(-6665, -7183, -5720)
    63 = @(-6788, -13220, -13162)
    64 = @(0, 0, 0)
    65 = @(0, 0, 0)
    66 = @(-7903, 5, -51)
    67 = @(0, 0, 0)
    68 = @(8048, 9613, 12877)
    69 = @(8169, 9714, 12943)
    70 = @(-121, -101, -66)
    74 = @(-1777, 3400, 7223)
    75 = @(0, 0, 0)
    76 = @(0, 0, 0)
    80 = @(27, 2907, 6551)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("N$row").Value = $vals[0]
    $ws.Range("O$row").Value = $vals[1]
    $ws.Range("P$row").Value = $vals[2]
}

# --- Row 79 (Part. de Acionistas Não Controladores): N/O carry values,
#     but P stays blank, matching the source data for this period. ---
$ws.Range("N79").Value = 0
$ws.Range("O79").Value = 0
$ws.Range("M79").Copy()
$ws.Range("P79").PasteSpecial(-4122)   # xlPasteFormats -> blank, untyped cell

# --- Section-header / subtotal-label rows: columns B:M are blank on these
#     rows, so N:P should likewise stay blank (no value/type), not zero. ---
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($row in $blankRows) {
    $ws.Range("M$row").Copy()
    $ws.Range("N${row}:P$row").PasteSpecial(-4122)   # xlPasteFormats
}
